# Daily scrape update - 2026-01-20 03:38:26 UTC
# Replaces the 3 sample opportunity rows with the 7 freshly scraped rows,
# highlights the "PREMIUM" = Yes cells, and resizes a few columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New opportunity rows (row 2 .. row 8), columns A-H:
#   A OPPORTUNITY ID, B OPPORTUNITY LINK, C TITLE, D COUNTRY,
#   E PREMIUM, F APPLICANTS, G DURATION, H ORGANIZATION
# ---------------------------------------------------------------------
$rows = @(
  @{ id="1331247"; title="Supply Planner Trainee";                        country="Panamá, Provincia de Panamá, Panamá"; premium="No";  applicants="3 applicants";   duration="6 - 18 Months"; org="NESTLE" }
  @{ id="1331240"; title="Sustainability Project Management trainee";     country="Bruxelles, Belgio";                    premium="No";  applicants="6 applicants";   duration="6 - 18 Months"; org="UCB" }
  @{ id="1331236"; title="[EXP] Routing & Capacity Analyst (EU Only)";    country="Bruxelles, Belgium";                   premium="Yes"; applicants="8 applicants";   duration="6 - 18 Months"; org="DHL Group" }
  @{ id="1330717"; title="Business Development Executive";                country="Rawalpindi, Pakistan";                 premium="No";  applicants="3 applicants";   duration="9 - 12 Weeks";  org="GrowUp Tech Solution" }
  @{ id="1330625"; title="Marketing Executive";                           country="Islamabad, Pakistan";                  premium="No";  applicants="1 applicant";    duration="9 - 12 Weeks";  org="GrowUp Tech Solution" }
  @{ id="1330623"; title="Intern – Strategy and Planning";                country="Karachi, Pakistan";                    premium="No";  applicants="8 applicants";   duration="9 - 12 Weeks";  org="Keys Productions (Pvt) Ltd" }
  @{ id="1284813"; title="Junior Technical Support Engineer";             country="Budapest, Hungary";                    premium="Yes"; applicants="148 applicants"; duration="6 - 18 Months"; org="EATON" }
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row.id
  $ws.Cells.Item($r, 2).Value = "https://aiesec.org/opportunity/global-talent/" + $row.id
  $ws.Cells.Item($r, 3).Value = $row.title
  $ws.Cells.Item($r, 4).Value = $row.country
  $ws.Cells.Item($r, 5).Value = $row.premium
  $ws.Cells.Item($r, 6).Value = $row.applicants
  $ws.Cells.Item($r, 7).Value = $row.duration
  $ws.Cells.Item($r, 8).Value = $row.org
  $r = $r + 1
}

# Highlight the PREMIUM = "Yes" cells (rows 4 and 8) with a yellow fill.
# (Applied as two separate single-area assignments: multi-area "E4,E8"
# unions only reliably format the first area on this host.)
$ws.Range("E4").Interior.ColorIndex = 6
$ws.Range("E8").Interior.ColorIndex = 6

# ---------------------------------------------------------------------
# Column width tweaks (C, D, G, H). Excel's COM layer stores column
# widths with a small fixed offset vs. the ColumnWidth value you set,
# so we back the offset out to land exactly on the target widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 44.166666666666664
$ws.Columns.Item(4).ColumnWidth = 37.166666666666664
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 28.166666666666668
